$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the time display: C4 used to hold "Female" (Gender) - it now shows a time string.
$ws.Range("C4").Value = "13:14:15 PM"

# The previous Gender value ("Female") and remaining row data move down to a new row 5,
# with the rest of row 5 filled with placeholder values.
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "r"
$ws.Range("C5").Value = "Female"
$ws.Range("D5").Value = "r"
$ws.Range("E5").Value = "06/07/2023"
$ws.Range("F5").Value = "r"
$ws.Range("G5").Value = "rr"
$ws.Range("H5").Value = "r"
